$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Update the cached "datetimeFigureOut" footer-date field from
#    09/06/2020 -> 12/06/2020 everywhere it appears (the slide master and
#    every slide layout each carry their own copy of the placeholder).
# ---------------------------------------------------------------------------
$oldDate = "09/06/2020"
$newDate = "12/06/2020"

function Update-DatePlaceholder {
    param($shapes)

    for ($j = 1; $j -le $shapes.Count; $j++) {
        $shp = $shapes.Item($j)

        $isDatePlaceholder = $false
        try {
            # ppPlaceholderDate == 16
            if ($shp.PlaceholderFormat.Type -eq 16) {
                $isDatePlaceholder = $true
            }
        } catch {
            $isDatePlaceholder = $false
        }

        if ($isDatePlaceholder -and $shp.HasTextFrame) {
            $tr = $shp.TextFrame.TextRange
            if ($tr.Text -eq $oldDate) {
                $tr.Text = $newDate
            }
        }
    }
}

# Slide master.
$master = $p.SlideMaster
Update-DatePlaceholder $master.Shapes

# Every custom layout hanging off the master.
for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    $layout = $master.CustomLayouts.Item($li)
    Update-DatePlaceholder $layout.Shapes
}

# ---------------------------------------------------------------------------
# 2) Slide 1: "specimen requirements" -> "specimen requested".
# ---------------------------------------------------------------------------
$oldLabel = "specimen requirements"
$newLabel = "specimen requested"

$slide = $p.Slides.Item(1)
for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
    $shp = $slide.Shapes.Item($i)
    if ($shp.HasTextFrame) {
        if ($shp.TextFrame.HasText) {
            if ($shp.TextFrame.TextRange.Text -eq $oldLabel) {
                $shp.TextFrame.TextRange.Text = $newLabel
            }
        }
    }
}
